# Helper files from data requests
# - Sort Table1 by the "Search" column (I) ascending (was sorted desc by "Last Login" / K)
# - Freeze panes at the header row, leaving active selection at L10

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$table = $ws.ListObjects.Item("Table1")

$sortCol = $ws.Range("I1:I195")

$table.Sort.SortFields.Clear()
$table.Sort.SortFields.Add($sortCol, 0, 1) | Out-Null
$table.Sort.Header = 1
$table.Sort.Apply()

$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

$ws.Range("L10").Select()
